$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$enDash = [char]0x2013

# 1. The 05/09 10:15-11:50 entry actually ran until 12:15, not 11:50 -> fix the end time.
$d.Content.Find.Execute(
    "05/09/2024, 10:15 $enDash 11:50", $false, $false, $false, $false, $false,
    $true, 1, $false, "05/09/2024, 10:15 $enDash 12:15", 2)

# 2. Fill in the two following (previously empty) log rows with the new JS news work.
$t.Cell(14, 1).Range.Text = "05/09/2024, 13:15 $enDash 14:00"
$t.Cell(14, 2).Range.Text = "Start gemaakt met nieuws via API"

$t.Cell(15, 1).Range.Text = "05/09/2024, 14:15 $enDash 16:00"
$t.Cell(15, 2).Range.Text = "Nieuws via fetch ingeladen, skeleton loaders voor nieuws gemaakt"
